$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-17
# from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175)
$newDate = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
